$d = $word.ActiveDocument
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false

function ReplaceText($findText, $isWildcard, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $isWildcard, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: $findText"
        return
    }
    $rng.Text = $newText
}

$quote = [char]34
$apos = [char]0x2019

# -----------------------------------------------------------------
# 1) Chief Justice ... crime. - merge split runs (text content identical)
# -----------------------------------------------------------------
ReplaceText "Chief Justice Arifin Zakaria dismissed*crime." $true "TMP_PLACEHOLDER_1"
$txt1 = "Chief Justice Arifin Zakaria dismissed Anwar's challenge against last year's conviction for sodomizing a young male former aide, saying the court found " + $quote + "overwhelming evidence" + $quote + " of the crime."
ReplaceText "TMP_PLACEHOLDER_1" $false $txt1

# -----------------------------------------------------------------
# 3) Analysts say jailing Anwar ... gerrymandering. - merge split runs
# -----------------------------------------------------------------
ReplaceText "Analysts say jailing Anwar also bears risks for*gerrymandering." $true "TMP_PLACEHOLDER_2"
$txt2 = "Analysts say jailing Anwar also bears risks for Najib's regime, which lost the popular vote in 2013 elections, clinging to power only through gerrymandering."
ReplaceText "TMP_PLACEHOLDER_2" $false $txt2

# -----------------------------------------------------------------
# 4) Remove UMNO text (row 12 col 2-4) -- do this BEFORE changing row 10's
#    UMNO text, since after that change the two cells would have identical
#    text and Find would match the wrong (first) occurrence.
# -----------------------------------------------------------------
ReplaceText "UMNO  has dominated reform discriminatory policies" $false ""

# -----------------------------------------------------------------
# 2) UMNO text change (row 10 col 2-4)
# -----------------------------------------------------------------
ReplaceText "UMNO  were filed opposition election showings" $false "UMNO  has dominated reform discriminatory policies"

# -----------------------------------------------------------------
# 5) Coh-Metrix Average - merge split runs
# -----------------------------------------------------------------
ReplaceText "Coh-Metrix Average" $false "TMP_PLACEHOLDER_3"
ReplaceText "TMP_PLACEHOLDER_3" $false "Coh-Metrix Average"

# -----------------------------------------------------------------
# 6) Flesch Kincaid Grade Level - merge split runs
# -----------------------------------------------------------------
ReplaceText "Flesch Kincaid Grade Level" $false "TMP_PLACEHOLDER_4"
ReplaceText "TMP_PLACEHOLDER_4" $false "Flesch Kincaid Grade Level"

# -----------------------------------------------------------------
# 7) What are the risks to Najib's regime for Anwar's conviction? - merge split runs
# -----------------------------------------------------------------
ReplaceText "What are the risks to*conviction?" $true "TMP_PLACEHOLDER_5"
$txt3 = "What are the risks to Najib" + $apos + "s regime for Anwar" + $apos + "s conviction?"
ReplaceText "TMP_PLACEHOLDER_5" $false $txt3

# -----------------------------------------------------------------
# 8) Numeric value insertions in previously-empty cells
# -----------------------------------------------------------------
$t = $d.Tables.Item(1)

$cell = $t.Cell(24, 4)
$cell.Range.Text = "1"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12

$cell = $t.Cell(44, 4)
$cell.Range.Text = "0.6666"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12

$cell = $t.Cell(45, 4)
$cell.Range.Text = "0.5"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12

$cell = $t.Cell(46, 4)
$cell.Range.Text = "0.5833"
$cell.Range.Font.Bold = 1
$cell.Range.Font.Size = 12
$cell.Range.Font.SizeBi = 12
